$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (want-to-go count) for two events
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1090   # 南宁·熊喵M动漫嘉年华【免费】
$wsExpo.Range("F4").Value = 2510   # 南宁·第二届北极光动漫展

# Sheet "全部类型" (all types): same two events appear again at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1090    # 南宁·熊喵M动漫嘉年华【免费】
$wsAll.Range("F6").Value = 2510    # 南宁·第二届北极光动漫展
